$d = $word.ActiveDocument
try {
  $full = $d.Content.XML()
  Write-Host ($full.Length)
} catch { Write-Host "ERR $_" }
